$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.663.24"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.612.80"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "513.02"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.79"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "2.628.15"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.80"
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "3.072.75"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "60.584.41"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.62"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "2.619.52"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.74"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.63"
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.57"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.87"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D28").Value = "0.0₃0841"
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.41"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.85"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.57"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  -2.56%  "
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.867"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.48"
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.39"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.841"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.31"
$ws.Range("E42").Value = "  -7.32%  "
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.625"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0553"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.78"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "2.000.83"
$ws.Range("E51").Value = "  -3.07%  "
